# Update Name of Algo
# Applies the updated RandomForest-imputed values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.2155
$ws.Range("B3").Value = 5.647600000000001
$ws.Range("E3").Value = 16.3695
$ws.Range("E12").Value = 17.53970000000001
$ws.Range("B14").Value = 5.367600000000003
$ws.Range("B16").Value = 6.489500000000001
$ws.Range("C18").Value = -11.4815
$ws.Range("B21").Value = 9.3385
$ws.Range("B23").Value = 9.078000000000003
$ws.Range("C24").Value = -12.76449999999999
$ws.Range("E24").Value = 16.4641
$ws.Range("B25").Value = 5.665199999999998
$ws.Range("C25").Value = -11.619
$ws.Range("E25").Value = 16.92660000000001
$ws.Range("B26").Value = 6.346000000000006
$ws.Range("C27").Value = -12.46039999999999
$ws.Range("B29").Value = 5.095600000000001
$ws.Range("C30").Value = -13.307
$ws.Range("C31").Value = -13.40319999999999
$ws.Range("C39").Value = -12.43140000000001
$ws.Range("B40").Value = 8.976699999999999
$ws.Range("E41").Value = 16.63389999999999
$ws.Range("C42").Value = -12.2954
$ws.Range("C48").Value = -11.41969999999999
$ws.Range("E50").Value = 16.2284
$ws.Range("C51").Value = -11.35999999999999
$ws.Range("C52").Value = -11.132
$ws.Range("B53").Value = 5.389399999999998
$ws.Range("E53").Value = 16.51790000000001
$ws.Range("C55").Value = -13.8202
$ws.Range("C56").Value = -13.44959999999999
$ws.Range("E56").Value = 16.25630000000001
$ws.Range("B57").Value = 5.232599999999995
$ws.Range("C57").Value = -13.43099999999999
$ws.Range("E57").Value = 16.7435
$ws.Range("E58").Value = 16.28920000000002
$ws.Range("B59").Value = 4.5847
$ws.Range("C60").Value = -13.83229999999999
$ws.Range("E61").Value = 16.4659
$ws.Range("E63").Value = 17.44870000000002
$ws.Range("E64").Value = 17.5558
$ws.Range("B65").Value = 6.021400000000001
$ws.Range("B69").Value = 5.348299999999995
$ws.Range("E70").Value = 17.44710000000001
$ws.Range("E72").Value = 16.91090000000001
$ws.Range("C73").Value = -12.6986
$ws.Range("C74").Value = -12.542
$ws.Range("B79").Value = 9.230400000000005
$ws.Range("B83").Value = 5.2627
$ws.Range("E86").Value = 16.5
$ws.Range("C89").Value = -10.3672
$ws.Range("E89").Value = 17.61180000000003
$ws.Range("C90").Value = -12.5246
$ws.Range("B91").Value = 5.084700000000002
$ws.Range("C92").Value = -10.5774
$ws.Range("B93").Value = 5.763300000000002
$ws.Range("E98").Value = 15.8914
$ws.Range("B100").Value = 5.255399999999998
$ws.Range("E100").Value = 16.49300000000001
$ws.Range("E102").Value = 16.57359999999998
